$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprint "Documentation de Projet et shémas de la databse" block (rows 22-27):
# Adjust the existing total-hours entry and fill in the two previously
# empty task rows with the new work that was done plus their hours.
$ws.Range("D22").Value = 6.5

$ws.Range("B23").Value = "Création des dockerfiles ajout dans le docker compose et créations des routes du backend"
$ws.Range("D23").Value = 2

$ws.Range("B24").Value = "Documentation et tests"
$ws.Range("D24").Value = 3.5

# Personal reflection for that sprint (row 28, merged B28:D28)
$ws.Range("B28").Value = "Ce sprint a été très dure un début avec beaucoup de tests non concluant ce qui n'était pas très motivant cel a été emplifier du fait que j'était malade mais une fois la solution trouvé j'ai pu m'attaquer a des choses plus ""tranquille"" ce qui a été mieux"

# Reflect where the author ended up looking while wrapping up this entry.
$ws.Range("E26").Select()
